$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.75%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "48.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "10.64%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.298"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.12%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07922"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.84%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.578"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.331"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "24.46%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.625"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.05%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1247"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.07%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1965"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.84%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09589"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.80%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04544"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "8.74%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1049"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.27%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001307"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.75%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04213"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.05%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005904"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.69%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.344"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.66%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.435"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.43%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.32%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.093"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.09%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1404"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.47%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.25%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001299"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.86%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.70%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001358"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.97%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003557"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.20%"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "0.06%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05858"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.83%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01038"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "85.50%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008013"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.75%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1449"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.28%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007554"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.22%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007932"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.50%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3186"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.50%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007023"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.05%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.91%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05588"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.21%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004020"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "1.93%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002110"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.91%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002010"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "1.91%"
